# This script applies the data refresh for the "living_rooms" memory-task input file.
# The underlying study regenerated the per-subject stimulus list (block_total 3 -> 6,
# trial_total renumbered, new stimuli/img_*.png assignments, refreshed norming stats,
# and the single catch-trial row moved from row 15 to row 29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 6
$ws.Cells.Item(2, 6).Value = 284
$ws.Cells.Item(2, 9).Value = $null
$ws.Cells.Item(2, 10).Value = 'new'
$ws.Cells.Item(2, 11).Value = 'f'
$ws.Cells.Item(2, 12).Value = 'stimuli/img_b21d7.png'
$ws.Cells.Item(2, 13).Value = 27.75555555555556
$ws.Cells.Item(2, 14).Value = 13.86666666666667
$ws.Cells.Item(2, 15).Value = 20.81111111111111
$ws.Cells.Item(2, 16).Value = 45
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = 1
$ws.Cells.Item(2, 19).Value = 1
$ws.Cells.Item(2, 20).Value = 2
$ws.Cells.Item(2, 21).Value = 2
$ws.Cells.Item(2, 22).Value = 1

# Row 3
$ws.Cells.Item(3, 3).Value = 6
$ws.Cells.Item(3, 6).Value = 285
$ws.Cells.Item(3, 12).Value = 'stimuli/img_tn8ys.png'
$ws.Cells.Item(3, 13).Value = 86.70454545454545
$ws.Cells.Item(3, 14).Value = 72.4090909090909
$ws.Cells.Item(3, 15).Value = 79.55681818181819
$ws.Cells.Item(3, 16).Value = 44
$ws.Cells.Item(3, 17).Value = 10
$ws.Cells.Item(3, 18).Value = 10
$ws.Cells.Item(3, 19).Value = 10
$ws.Cells.Item(3, 20).Value = 9
$ws.Cells.Item(3, 21).Value = 9
$ws.Cells.Item(3, 22).Value = 10

# Row 4
$ws.Cells.Item(4, 3).Value = 6
$ws.Cells.Item(4, 6).Value = 286
$ws.Cells.Item(4, 9).Value = $null
$ws.Cells.Item(4, 10).Value = 'new'
$ws.Cells.Item(4, 11).Value = 'f'
$ws.Cells.Item(4, 12).Value = 'stimuli/img_x4bln.png'
$ws.Cells.Item(4, 13).Value = 76.34042553191489
$ws.Cells.Item(4, 14).Value = 59.51063829787234
$ws.Cells.Item(4, 15).Value = 67.92553191489361
$ws.Cells.Item(4, 16).Value = 47
$ws.Cells.Item(4, 17).Value = 7
$ws.Cells.Item(4, 18).Value = 7
$ws.Cells.Item(4, 19).Value = 7
$ws.Cells.Item(4, 20).Value = 7
$ws.Cells.Item(4, 21).Value = 7
$ws.Cells.Item(4, 22).Value = 7

# Row 5
$ws.Cells.Item(5, 3).Value = 6
$ws.Cells.Item(5, 6).Value = 287
$ws.Cells.Item(5, 9).Value = 'target'
$ws.Cells.Item(5, 10).Value = 'old'
$ws.Cells.Item(5, 11).Value = 'j'
$ws.Cells.Item(5, 12).Value = 'stimuli/img_a9he3.png'
$ws.Cells.Item(5, 13).Value = 83.06521739130434
$ws.Cells.Item(5, 14).Value = 63.95652173913044
$ws.Cells.Item(5, 15).Value = 73.51086956521739
$ws.Cells.Item(5, 16).Value = 46
$ws.Cells.Item(5, 17).Value = 8
$ws.Cells.Item(5, 18).Value = 8
$ws.Cells.Item(5, 19).Value = 8
$ws.Cells.Item(5, 20).Value = 8
$ws.Cells.Item(5, 21).Value = 8
$ws.Cells.Item(5, 22).Value = 8

# Row 6
$ws.Cells.Item(6, 3).Value = 6
$ws.Cells.Item(6, 6).Value = 288
$ws.Cells.Item(6, 9).Value = $null
$ws.Cells.Item(6, 10).Value = 'new'
$ws.Cells.Item(6, 11).Value = 'f'
$ws.Cells.Item(6, 12).Value = 'stimuli/img_c89x3.png'
$ws.Cells.Item(6, 13).Value = 72.8695652173913
$ws.Cells.Item(6, 14).Value = 49.65217391304348
$ws.Cells.Item(6, 15).Value = 61.26086956521739
$ws.Cells.Item(6, 17).Value = 6
$ws.Cells.Item(6, 18).Value = 6
$ws.Cells.Item(6, 19).Value = 6
$ws.Cells.Item(6, 20).Value = 6
$ws.Cells.Item(6, 21).Value = 6
$ws.Cells.Item(6, 22).Value = 5

# Row 7
$ws.Cells.Item(7, 3).Value = 6
$ws.Cells.Item(7, 6).Value = 289
$ws.Cells.Item(7, 12).Value = 'stimuli/img_o30wb.png'
$ws.Cells.Item(7, 13).Value = 81.06666666666666
$ws.Cells.Item(7, 14).Value = 65.37777777777778
$ws.Cells.Item(7, 15).Value = 73.22222222222223
$ws.Cells.Item(7, 16).Value = 45
$ws.Cells.Item(7, 17).Value = 8
$ws.Cells.Item(7, 18).Value = 8
$ws.Cells.Item(7, 19).Value = 8
$ws.Cells.Item(7, 20).Value = 8
$ws.Cells.Item(7, 21).Value = 8
$ws.Cells.Item(7, 22).Value = 8

# Row 8
$ws.Cells.Item(8, 3).Value = 6
$ws.Cells.Item(8, 6).Value = 290
$ws.Cells.Item(8, 9).Value = $null
$ws.Cells.Item(8, 10).Value = 'new'
$ws.Cells.Item(8, 11).Value = 'f'
$ws.Cells.Item(8, 12).Value = 'stimuli/img_vgh2g.png'
$ws.Cells.Item(8, 13).Value = 93.81395348837209
$ws.Cells.Item(8, 14).Value = 78.27906976744185
$ws.Cells.Item(8, 15).Value = 86.04651162790697
$ws.Cells.Item(8, 16).Value = 43
$ws.Cells.Item(8, 17).Value = 10
$ws.Cells.Item(8, 18).Value = 10
$ws.Cells.Item(8, 19).Value = 10
$ws.Cells.Item(8, 20).Value = 10
$ws.Cells.Item(8, 21).Value = 10
$ws.Cells.Item(8, 22).Value = 10

# Row 9
$ws.Cells.Item(9, 3).Value = 6
$ws.Cells.Item(9, 6).Value = 291
$ws.Cells.Item(9, 9).Value = 'target'
$ws.Cells.Item(9, 10).Value = 'old'
$ws.Cells.Item(9, 11).Value = 'j'
$ws.Cells.Item(9, 12).Value = 'stimuli/img_ra2nm.png'
$ws.Cells.Item(9, 13).Value = 70.75
$ws.Cells.Item(9, 14).Value = 50.375
$ws.Cells.Item(9, 15).Value = 60.5625
$ws.Cells.Item(9, 16).Value = 40
$ws.Cells.Item(9, 20).Value = 5
$ws.Cells.Item(9, 21).Value = 5

# Row 10
$ws.Cells.Item(10, 3).Value = 6
$ws.Cells.Item(10, 6).Value = 292
$ws.Cells.Item(10, 12).Value = 'stimuli/img_zxvl3.png'
$ws.Cells.Item(10, 13).Value = 68.78260869565217
$ws.Cells.Item(10, 14).Value = 47.56521739130435
$ws.Cells.Item(10, 15).Value = 58.17391304347827
$ws.Cells.Item(10, 16).Value = 46
$ws.Cells.Item(10, 17).Value = 5
$ws.Cells.Item(10, 18).Value = 5
$ws.Cells.Item(10, 19).Value = 5
$ws.Cells.Item(10, 20).Value = 5
$ws.Cells.Item(10, 21).Value = 5
$ws.Cells.Item(10, 22).Value = 5

# Row 11
$ws.Cells.Item(11, 3).Value = 6
$ws.Cells.Item(11, 6).Value = 293
$ws.Cells.Item(11, 12).Value = 'stimuli/img_z4jxm.png'
$ws.Cells.Item(11, 13).Value = 88.30952380952381
$ws.Cells.Item(11, 14).Value = 72.64285714285714
$ws.Cells.Item(11, 15).Value = 80.47619047619048
$ws.Cells.Item(11, 16).Value = 42
$ws.Cells.Item(11, 17).Value = 10
$ws.Cells.Item(11, 18).Value = 10
$ws.Cells.Item(11, 19).Value = 10
$ws.Cells.Item(11, 20).Value = 10
$ws.Cells.Item(11, 21).Value = 10
$ws.Cells.Item(11, 22).Value = 10

# Row 12
$ws.Cells.Item(12, 3).Value = 6
$ws.Cells.Item(12, 6).Value = 294
$ws.Cells.Item(12, 9).Value = 'target'
$ws.Cells.Item(12, 10).Value = 'old'
$ws.Cells.Item(12, 11).Value = 'j'
$ws.Cells.Item(12, 12).Value = 'stimuli/img_5jy9c.png'
$ws.Cells.Item(12, 13).Value = 87.37209302325581
$ws.Cells.Item(12, 14).Value = 79.18604651162791
$ws.Cells.Item(12, 15).Value = 83.27906976744185
$ws.Cells.Item(12, 17).Value = 10
$ws.Cells.Item(12, 18).Value = 10
$ws.Cells.Item(12, 19).Value = 10
$ws.Cells.Item(12, 20).Value = 10
$ws.Cells.Item(12, 21).Value = 9
$ws.Cells.Item(12, 22).Value = 10

# Row 13
$ws.Cells.Item(13, 3).Value = 6
$ws.Cells.Item(13, 6).Value = 295
$ws.Cells.Item(13, 9).Value = 'target'
$ws.Cells.Item(13, 10).Value = 'old'
$ws.Cells.Item(13, 11).Value = 'j'
$ws.Cells.Item(13, 12).Value = 'stimuli/img_9bkl9.png'
$ws.Cells.Item(13, 13).Value = 46.62162162162162
$ws.Cells.Item(13, 14).Value = 34.27027027027027
$ws.Cells.Item(13, 15).Value = 40.44594594594595
$ws.Cells.Item(13, 16).Value = 37
$ws.Cells.Item(13, 20).Value = 3
$ws.Cells.Item(13, 22).Value = 3

# Row 14
$ws.Cells.Item(14, 3).Value = 6
$ws.Cells.Item(14, 6).Value = 296
$ws.Cells.Item(14, 9).Value = $null
$ws.Cells.Item(14, 10).Value = 'new'
$ws.Cells.Item(14, 11).Value = 'f'
$ws.Cells.Item(14, 12).Value = 'stimuli/img_pdzf1.png'
$ws.Cells.Item(14, 13).Value = 86.23913043478261
$ws.Cells.Item(14, 14).Value = 67.17391304347827
$ws.Cells.Item(14, 15).Value = 76.70652173913044
$ws.Cells.Item(14, 17).Value = 9
$ws.Cells.Item(14, 18).Value = 9
$ws.Cells.Item(14, 19).Value = 9
$ws.Cells.Item(14, 20).Value = 9
$ws.Cells.Item(14, 21).Value = 9
$ws.Cells.Item(14, 22).Value = 8

# Row 15
$ws.Cells.Item(15, 3).Value = 6
$ws.Cells.Item(15, 6).Value = 297
$ws.Cells.Item(15, 8).Value = 'living_rooms'
$ws.Cells.Item(15, 9).Value = 'target'
$ws.Cells.Item(15, 10).Value = 'old'
$ws.Cells.Item(15, 11).Value = 'j'
$ws.Cells.Item(15, 12).Value = 'stimuli/img_g13d5.png'
$ws.Cells.Item(15, 13).Value = 73
$ws.Cells.Item(15, 14).Value = 51.51111111111111
$ws.Cells.Item(15, 15).Value = 62.25555555555556
$ws.Cells.Item(15, 16).Value = 45
$ws.Cells.Item(15, 17).Value = 6
$ws.Cells.Item(15, 18).Value = 6
$ws.Cells.Item(15, 19).Value = 6
$ws.Cells.Item(15, 20).Value = 6
$ws.Cells.Item(15, 21).Value = 6
$ws.Cells.Item(15, 22).Value = 6

# Row 16
$ws.Cells.Item(16, 3).Value = 6
$ws.Cells.Item(16, 6).Value = 298
$ws.Cells.Item(16, 9).Value = $null
$ws.Cells.Item(16, 10).Value = 'new'
$ws.Cells.Item(16, 11).Value = 'f'
$ws.Cells.Item(16, 12).Value = 'stimuli/img_xr3up.png'
$ws.Cells.Item(16, 13).Value = 76.24444444444444
$ws.Cells.Item(16, 14).Value = 55.88888888888889
$ws.Cells.Item(16, 15).Value = 66.06666666666666
$ws.Cells.Item(16, 16).Value = 45
$ws.Cells.Item(16, 17).Value = 7
$ws.Cells.Item(16, 18).Value = 7
$ws.Cells.Item(16, 19).Value = 7
$ws.Cells.Item(16, 20).Value = 6
$ws.Cells.Item(16, 21).Value = 6
$ws.Cells.Item(16, 22).Value = 6

# Row 17
$ws.Cells.Item(17, 3).Value = 6
$ws.Cells.Item(17, 6).Value = 299
$ws.Cells.Item(17, 12).Value = 'stimuli/img_165pk.png'
$ws.Cells.Item(17, 13).Value = 85.73333333333333
$ws.Cells.Item(17, 14).Value = 69.22222222222223
$ws.Cells.Item(17, 15).Value = 77.47777777777779
$ws.Cells.Item(17, 16).Value = 45
$ws.Cells.Item(17, 17).Value = 9
$ws.Cells.Item(17, 18).Value = 9
$ws.Cells.Item(17, 19).Value = 9
$ws.Cells.Item(17, 20).Value = 9
$ws.Cells.Item(17, 21).Value = 9
$ws.Cells.Item(17, 22).Value = 9

# Row 18
$ws.Cells.Item(18, 3).Value = 6
$ws.Cells.Item(18, 6).Value = 300
$ws.Cells.Item(18, 9).Value = $null
$ws.Cells.Item(18, 10).Value = 'new'
$ws.Cells.Item(18, 11).Value = 'f'
$ws.Cells.Item(18, 12).Value = 'stimuli/img_wbws6.png'
$ws.Cells.Item(18, 13).Value = 57.97777777777777
$ws.Cells.Item(18, 14).Value = 42.53333333333333
$ws.Cells.Item(18, 15).Value = 50.25555555555555
$ws.Cells.Item(18, 16).Value = 45
$ws.Cells.Item(18, 17).Value = 4
$ws.Cells.Item(18, 18).Value = 4
$ws.Cells.Item(18, 19).Value = 4
$ws.Cells.Item(18, 20).Value = 4
$ws.Cells.Item(18, 21).Value = 4

# Row 19
$ws.Cells.Item(19, 3).Value = 6
$ws.Cells.Item(19, 6).Value = 301
$ws.Cells.Item(19, 9).Value = 'target'
$ws.Cells.Item(19, 10).Value = 'old'
$ws.Cells.Item(19, 11).Value = 'j'
$ws.Cells.Item(19, 12).Value = 'stimuli/img_jpldg.png'
$ws.Cells.Item(19, 13).Value = 79.54545454545455
$ws.Cells.Item(19, 14).Value = 57.75
$ws.Cells.Item(19, 15).Value = 68.64772727272728
$ws.Cells.Item(19, 16).Value = 44
$ws.Cells.Item(19, 17).Value = 7
$ws.Cells.Item(19, 18).Value = 7
$ws.Cells.Item(19, 19).Value = 7
$ws.Cells.Item(19, 20).Value = 7
$ws.Cells.Item(19, 21).Value = 7
$ws.Cells.Item(19, 22).Value = 7

# Row 20
$ws.Cells.Item(20, 3).Value = 6
$ws.Cells.Item(20, 6).Value = 302
$ws.Cells.Item(20, 12).Value = 'stimuli/img_j4ttn.png'
$ws.Cells.Item(20, 13).Value = 12.61904761904762
$ws.Cells.Item(20, 14).Value = 11.42857142857143
$ws.Cells.Item(20, 15).Value = 12.02380952380952
$ws.Cells.Item(20, 16).Value = 42
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = 1
$ws.Cells.Item(20, 19).Value = 1
$ws.Cells.Item(20, 20).Value = 1
$ws.Cells.Item(20, 21).Value = 1
$ws.Cells.Item(20, 22).Value = 1

# Row 21
$ws.Cells.Item(21, 3).Value = 6
$ws.Cells.Item(21, 6).Value = 303
$ws.Cells.Item(21, 9).Value = 'target'
$ws.Cells.Item(21, 10).Value = 'old'
$ws.Cells.Item(21, 11).Value = 'j'
$ws.Cells.Item(21, 12).Value = 'stimuli/img_3sw8t.png'
$ws.Cells.Item(21, 13).Value = 67.4888888888889
$ws.Cells.Item(21, 14).Value = 48.51111111111111
$ws.Cells.Item(21, 15).Value = 58
$ws.Cells.Item(21, 16).Value = 45
$ws.Cells.Item(21, 17).Value = 5
$ws.Cells.Item(21, 18).Value = 5
$ws.Cells.Item(21, 19).Value = 5
$ws.Cells.Item(21, 20).Value = 5
$ws.Cells.Item(21, 21).Value = 5
$ws.Cells.Item(21, 22).Value = 5

# Row 22
$ws.Cells.Item(22, 3).Value = 6
$ws.Cells.Item(22, 6).Value = 304
$ws.Cells.Item(22, 12).Value = 'stimuli/img_pjfx6.png'
$ws.Cells.Item(22, 13).Value = 32.23404255319149
$ws.Cells.Item(22, 14).Value = 26.59574468085106
$ws.Cells.Item(22, 15).Value = 29.41489361702127
$ws.Cells.Item(22, 16).Value = 47
$ws.Cells.Item(22, 17).Value = 2
$ws.Cells.Item(22, 18).Value = 2
$ws.Cells.Item(22, 19).Value = 2
$ws.Cells.Item(22, 20).Value = 2
$ws.Cells.Item(22, 21).Value = 2
$ws.Cells.Item(22, 22).Value = 3

# Row 23
$ws.Cells.Item(23, 3).Value = 6
$ws.Cells.Item(23, 6).Value = 305
$ws.Cells.Item(23, 9).Value = 'target'
$ws.Cells.Item(23, 10).Value = 'old'
$ws.Cells.Item(23, 11).Value = 'j'
$ws.Cells.Item(23, 12).Value = 'stimuli/img_pbsj1.png'
$ws.Cells.Item(23, 13).Value = 73.88636363636364
$ws.Cells.Item(23, 14).Value = 51.52272727272727
$ws.Cells.Item(23, 15).Value = 62.70454545454545
$ws.Cells.Item(23, 16).Value = 44
$ws.Cells.Item(23, 17).Value = 6
$ws.Cells.Item(23, 18).Value = 6
$ws.Cells.Item(23, 19).Value = 6
$ws.Cells.Item(23, 20).Value = 6
$ws.Cells.Item(23, 21).Value = 6
$ws.Cells.Item(23, 22).Value = 6

# Row 24
$ws.Cells.Item(24, 3).Value = 6
$ws.Cells.Item(24, 6).Value = 306
$ws.Cells.Item(24, 9).Value = $null
$ws.Cells.Item(24, 10).Value = 'new'
$ws.Cells.Item(24, 11).Value = 'f'
$ws.Cells.Item(24, 12).Value = 'stimuli/img_i6wsx.png'
$ws.Cells.Item(24, 13).Value = 79.07142857142857
$ws.Cells.Item(24, 14).Value = 58
$ws.Cells.Item(24, 15).Value = 68.53571428571428
$ws.Cells.Item(24, 17).Value = 7
$ws.Cells.Item(24, 18).Value = 7
$ws.Cells.Item(24, 19).Value = 7
$ws.Cells.Item(24, 20).Value = 7
$ws.Cells.Item(24, 21).Value = 7
$ws.Cells.Item(24, 22).Value = 7

# Row 25
$ws.Cells.Item(25, 3).Value = 6
$ws.Cells.Item(25, 6).Value = 307
$ws.Cells.Item(25, 9).Value = $null
$ws.Cells.Item(25, 10).Value = 'new'
$ws.Cells.Item(25, 11).Value = 'f'
$ws.Cells.Item(25, 12).Value = 'stimuli/img_xzyzy.png'
$ws.Cells.Item(25, 13).Value = 85.37209302325581
$ws.Cells.Item(25, 14).Value = 68.90697674418605
$ws.Cells.Item(25, 15).Value = 77.13953488372093
$ws.Cells.Item(25, 16).Value = 43
$ws.Cells.Item(25, 17).Value = 9
$ws.Cells.Item(25, 18).Value = 9
$ws.Cells.Item(25, 19).Value = 9
$ws.Cells.Item(25, 20).Value = 9
$ws.Cells.Item(25, 21).Value = 9
$ws.Cells.Item(25, 22).Value = 9

# Row 26
$ws.Cells.Item(26, 3).Value = 6
$ws.Cells.Item(26, 6).Value = 308
$ws.Cells.Item(26, 9).Value = $null
$ws.Cells.Item(26, 10).Value = 'new'
$ws.Cells.Item(26, 11).Value = 'f'
$ws.Cells.Item(26, 12).Value = 'stimuli/img_q9lab.png'
$ws.Cells.Item(26, 13).Value = 53.97560975609756
$ws.Cells.Item(26, 14).Value = 32.90243902439025
$ws.Cells.Item(26, 15).Value = 43.4390243902439
$ws.Cells.Item(26, 16).Value = 41
$ws.Cells.Item(26, 17).Value = 3
$ws.Cells.Item(26, 18).Value = 3
$ws.Cells.Item(26, 19).Value = 3
$ws.Cells.Item(26, 20).Value = 3
$ws.Cells.Item(26, 21).Value = 4
$ws.Cells.Item(26, 22).Value = 3

# Row 27
$ws.Cells.Item(27, 3).Value = 6
$ws.Cells.Item(27, 6).Value = 309
$ws.Cells.Item(27, 12).Value = 'stimuli/img_gka64.png'
$ws.Cells.Item(27, 13).Value = 19.23809523809524
$ws.Cells.Item(27, 14).Value = 20.02380952380953
$ws.Cells.Item(27, 15).Value = 19.63095238095238
$ws.Cells.Item(27, 17).Value = 1
$ws.Cells.Item(27, 18).Value = 1
$ws.Cells.Item(27, 19).Value = 1
$ws.Cells.Item(27, 20).Value = 1
$ws.Cells.Item(27, 21).Value = 1
$ws.Cells.Item(27, 22).Value = 2

# Row 28
$ws.Cells.Item(28, 3).Value = 6
$ws.Cells.Item(28, 6).Value = 310
$ws.Cells.Item(28, 9).Value = $null
$ws.Cells.Item(28, 10).Value = 'new'
$ws.Cells.Item(28, 11).Value = 'f'
$ws.Cells.Item(28, 12).Value = 'stimuli/img_c0vzo.png'
$ws.Cells.Item(28, 13).Value = 21.51162790697675
$ws.Cells.Item(28, 14).Value = 8.232558139534884
$ws.Cells.Item(28, 15).Value = 14.87209302325581
$ws.Cells.Item(28, 17).Value = 1
$ws.Cells.Item(28, 18).Value = 1
$ws.Cells.Item(28, 19).Value = 1
$ws.Cells.Item(28, 20).Value = 1
$ws.Cells.Item(28, 21).Value = 1
$ws.Cells.Item(28, 22).Value = 1

# Row 29
$ws.Cells.Item(29, 3).Value = 6
$ws.Cells.Item(29, 6).Value = 311
$ws.Cells.Item(29, 8).Value = $null
$ws.Cells.Item(29, 9).Value = $null
$ws.Cells.Item(29, 10).Value = 'catch'
$ws.Cells.Item(29, 11).Value = 'f'
$ws.Cells.Item(29, 12).Value = 'stimuli/catch_05.jpg'
$ws.Cells.Item(29, 13).Value = $null
$ws.Cells.Item(29, 14).Value = $null
$ws.Cells.Item(29, 15).Value = $null
$ws.Cells.Item(29, 16).Value = $null
$ws.Cells.Item(29, 17).Value = $null
$ws.Cells.Item(29, 18).Value = $null
$ws.Cells.Item(29, 19).Value = $null
$ws.Cells.Item(29, 20).Value = $null
$ws.Cells.Item(29, 21).Value = $null
$ws.Cells.Item(29, 22).Value = $null

# Row 30
$ws.Cells.Item(30, 3).Value = 6
$ws.Cells.Item(30, 6).Value = 312
$ws.Cells.Item(30, 9).Value = 'target'
$ws.Cells.Item(30, 10).Value = 'old'
$ws.Cells.Item(30, 11).Value = 'j'
$ws.Cells.Item(30, 12).Value = 'stimuli/img_wgddx.png'
$ws.Cells.Item(30, 13).Value = 45.6304347826087
$ws.Cells.Item(30, 14).Value = 34.30434782608695
$ws.Cells.Item(30, 15).Value = 39.96739130434783
$ws.Cells.Item(30, 16).Value = 46
$ws.Cells.Item(30, 17).Value = 3
$ws.Cells.Item(30, 18).Value = 3
$ws.Cells.Item(30, 19).Value = 3
$ws.Cells.Item(30, 20).Value = 3
$ws.Cells.Item(30, 21).Value = 3
$ws.Cells.Item(30, 22).Value = 4

# Row 31
$ws.Cells.Item(31, 3).Value = 6
$ws.Cells.Item(31, 6).Value = 313
$ws.Cells.Item(31, 12).Value = 'stimuli/img_16kib.png'
$ws.Cells.Item(31, 13).Value = 80.97727272727273
$ws.Cells.Item(31, 14).Value = 61.11363636363637
$ws.Cells.Item(31, 15).Value = 71.04545454545455
$ws.Cells.Item(31, 16).Value = 44
$ws.Cells.Item(31, 17).Value = 8
$ws.Cells.Item(31, 18).Value = 8
$ws.Cells.Item(31, 19).Value = 8
$ws.Cells.Item(31, 20).Value = 7
$ws.Cells.Item(31, 21).Value = 7
$ws.Cells.Item(31, 22).Value = 7

# Row 32
$ws.Cells.Item(32, 3).Value = 6
$ws.Cells.Item(32, 6).Value = 314
$ws.Cells.Item(32, 12).Value = 'stimuli/img_j856a.png'
$ws.Cells.Item(32, 13).Value = 38.225
$ws.Cells.Item(32, 14).Value = 25.875
$ws.Cells.Item(32, 15).Value = 32.05
$ws.Cells.Item(32, 16).Value = 40
$ws.Cells.Item(32, 17).Value = 2
$ws.Cells.Item(32, 18).Value = 2
$ws.Cells.Item(32, 19).Value = 2
$ws.Cells.Item(32, 20).Value = 3
$ws.Cells.Item(32, 21).Value = 3
$ws.Cells.Item(32, 22).Value = 2

# Row 33
$ws.Cells.Item(33, 3).Value = 6
$ws.Cells.Item(33, 6).Value = 315
$ws.Cells.Item(33, 9).Value = 'target'
$ws.Cells.Item(33, 10).Value = 'old'
$ws.Cells.Item(33, 11).Value = 'j'
$ws.Cells.Item(33, 12).Value = 'stimuli/img_w8yhd.png'
$ws.Cells.Item(33, 13).Value = 55.74418604651163
$ws.Cells.Item(33, 14).Value = 38.90697674418605
$ws.Cells.Item(33, 15).Value = 47.32558139534883
$ws.Cells.Item(33, 16).Value = 43
$ws.Cells.Item(33, 17).Value = 4
$ws.Cells.Item(33, 18).Value = 4
$ws.Cells.Item(33, 19).Value = 4
$ws.Cells.Item(33, 20).Value = 4
$ws.Cells.Item(33, 21).Value = 4
$ws.Cells.Item(33, 22).Value = 4

# Row 34
$ws.Cells.Item(34, 3).Value = 6
$ws.Cells.Item(34, 6).Value = 316
$ws.Cells.Item(34, 9).Value = $null
$ws.Cells.Item(34, 10).Value = 'new'
$ws.Cells.Item(34, 11).Value = 'f'
$ws.Cells.Item(34, 12).Value = 'stimuli/img_hmmra.png'
$ws.Cells.Item(34, 13).Value = 54.65853658536585
$ws.Cells.Item(34, 14).Value = 34.24390243902439
$ws.Cells.Item(34, 15).Value = 44.45121951219512
$ws.Cells.Item(34, 16).Value = 41
$ws.Cells.Item(34, 17).Value = 3
$ws.Cells.Item(34, 18).Value = 3
$ws.Cells.Item(34, 19).Value = 3
$ws.Cells.Item(34, 20).Value = 4
$ws.Cells.Item(34, 21).Value = 4
$ws.Cells.Item(34, 22).Value = 3

# Row 35
$ws.Cells.Item(35, 3).Value = 6
$ws.Cells.Item(35, 6).Value = 317
$ws.Cells.Item(35, 9).Value = 'target'
$ws.Cells.Item(35, 10).Value = 'old'
$ws.Cells.Item(35, 11).Value = 'j'
$ws.Cells.Item(35, 12).Value = 'stimuli/img_8dmpq.png'
$ws.Cells.Item(35, 13).Value = 30.65909090909091
$ws.Cells.Item(35, 14).Value = 24.11363636363636
$ws.Cells.Item(35, 15).Value = 27.38636363636364
$ws.Cells.Item(35, 17).Value = 2
$ws.Cells.Item(35, 18).Value = 2
$ws.Cells.Item(35, 19).Value = 2
$ws.Cells.Item(35, 20).Value = 2
$ws.Cells.Item(35, 21).Value = 2
$ws.Cells.Item(35, 22).Value = 2

# Row 36
$ws.Cells.Item(36, 3).Value = 6
$ws.Cells.Item(36, 6).Value = 318
$ws.Cells.Item(36, 12).Value = 'stimuli/img_6a0hu.png'
$ws.Cells.Item(36, 13).Value = 61.275
$ws.Cells.Item(36, 14).Value = 42.025
$ws.Cells.Item(36, 15).Value = 51.65
$ws.Cells.Item(36, 16).Value = 40
$ws.Cells.Item(36, 17).Value = 4
$ws.Cells.Item(36, 18).Value = 4
$ws.Cells.Item(36, 19).Value = 4
$ws.Cells.Item(36, 20).Value = 5
$ws.Cells.Item(36, 21).Value = 4
$ws.Cells.Item(36, 22).Value = 5

# Row 37
$ws.Cells.Item(37, 3).Value = 6
$ws.Cells.Item(37, 6).Value = 319
$ws.Cells.Item(37, 9).Value = 'target'
$ws.Cells.Item(37, 10).Value = 'old'
$ws.Cells.Item(37, 11).Value = 'j'
$ws.Cells.Item(37, 12).Value = 'stimuli/img_jkm86.png'
$ws.Cells.Item(37, 13).Value = 58.32558139534883
$ws.Cells.Item(37, 14).Value = 38.65116279069768
$ws.Cells.Item(37, 15).Value = 48.48837209302326
$ws.Cells.Item(37, 16).Value = 43
$ws.Cells.Item(37, 17).Value = 4
$ws.Cells.Item(37, 18).Value = 4
$ws.Cells.Item(37, 19).Value = 4
$ws.Cells.Item(37, 20).Value = 4
$ws.Cells.Item(37, 21).Value = 4
$ws.Cells.Item(37, 22).Value = 4

# Row 38
$ws.Cells.Item(38, 3).Value = 6
$ws.Cells.Item(38, 6).Value = 320
$ws.Cells.Item(38, 9).Value = $null
$ws.Cells.Item(38, 10).Value = 'new'
$ws.Cells.Item(38, 11).Value = 'f'
$ws.Cells.Item(38, 12).Value = 'stimuli/img_vgaye.png'
$ws.Cells.Item(38, 13).Value = 80.33333333333333
$ws.Cells.Item(38, 14).Value = 64.57777777777778
$ws.Cells.Item(38, 15).Value = 72.45555555555555
$ws.Cells.Item(38, 16).Value = 45
$ws.Cells.Item(38, 17).Value = 8
$ws.Cells.Item(38, 18).Value = 8
$ws.Cells.Item(38, 19).Value = 8
$ws.Cells.Item(38, 20).Value = 8
$ws.Cells.Item(38, 21).Value = 7
$ws.Cells.Item(38, 22).Value = 8

# Row 39
$ws.Cells.Item(39, 3).Value = 6
$ws.Cells.Item(39, 6).Value = 321
$ws.Cells.Item(39, 9).Value = $null
$ws.Cells.Item(39, 10).Value = 'new'
$ws.Cells.Item(39, 11).Value = 'f'
$ws.Cells.Item(39, 12).Value = 'stimuli/img_rych7.png'
$ws.Cells.Item(39, 13).Value = 30.4468085106383
$ws.Cells.Item(39, 14).Value = 23.4468085106383
$ws.Cells.Item(39, 15).Value = 26.9468085106383
$ws.Cells.Item(39, 17).Value = 2
$ws.Cells.Item(39, 18).Value = 2
$ws.Cells.Item(39, 19).Value = 2
$ws.Cells.Item(39, 20).Value = 2
$ws.Cells.Item(39, 21).Value = 2
$ws.Cells.Item(39, 22).Value = 2

# Row 40
$ws.Cells.Item(40, 3).Value = 6
$ws.Cells.Item(40, 6).Value = 322
$ws.Cells.Item(40, 9).Value = $null
$ws.Cells.Item(40, 10).Value = 'new'
$ws.Cells.Item(40, 11).Value = 'f'
$ws.Cells.Item(40, 12).Value = 'stimuli/img_nb8p4.png'
$ws.Cells.Item(40, 13).Value = 16.36170212765957
$ws.Cells.Item(40, 14).Value = 12.70212765957447
$ws.Cells.Item(40, 15).Value = 14.53191489361702
$ws.Cells.Item(40, 16).Value = 47
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = 1
$ws.Cells.Item(40, 19).Value = 1
$ws.Cells.Item(40, 20).Value = 1
$ws.Cells.Item(40, 21).Value = 1
$ws.Cells.Item(40, 22).Value = 1

# Row 41
$ws.Cells.Item(41, 3).Value = 6
$ws.Cells.Item(41, 6).Value = 323
$ws.Cells.Item(41, 12).Value = 'stimuli/img_5jp4f.png'
$ws.Cells.Item(41, 13).Value = 84.85714285714286
$ws.Cells.Item(41, 14).Value = 67.83333333333333
$ws.Cells.Item(41, 15).Value = 76.3452380952381
$ws.Cells.Item(41, 16).Value = 42
$ws.Cells.Item(41, 17).Value = 9
$ws.Cells.Item(41, 18).Value = 9
$ws.Cells.Item(41, 19).Value = 9
$ws.Cells.Item(41, 20).Value = 8
$ws.Cells.Item(41, 21).Value = 8
$ws.Cells.Item(41, 22).Value = 9

# Row 42
$ws.Cells.Item(42, 3).Value = 6
$ws.Cells.Item(42, 6).Value = 324
$ws.Cells.Item(42, 9).Value = $null
$ws.Cells.Item(42, 10).Value = 'new'
$ws.Cells.Item(42, 11).Value = 'f'
$ws.Cells.Item(42, 12).Value = 'stimuli/img_dg5h7.png'
$ws.Cells.Item(42, 13).Value = 88.72093023255815
$ws.Cells.Item(42, 14).Value = 76.06976744186046
$ws.Cells.Item(42, 15).Value = 82.3953488372093
$ws.Cells.Item(42, 16).Value = 43
$ws.Cells.Item(42, 17).Value = 10
$ws.Cells.Item(42, 18).Value = 10
$ws.Cells.Item(42, 19).Value = 10
$ws.Cells.Item(42, 20).Value = 10
$ws.Cells.Item(42, 21).Value = 10
$ws.Cells.Item(42, 22).Value = 10

